$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.860.09"
$ws.Range("E2").Value = "  -4.15%  "

# Row 3
$ws.Range("D3").Value = "2.453.33"
$ws.Range("E3").Value = "  -5.84%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'549.74"
$ws.Range("E5").Value = "  -3.76%  "

# Row 6
$ws.Range("D6").Value = "'145.22"
$ws.Range("E6").Value = "  -5.74%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.598"
$ws.Range("E8").Value = "  -3.33%  "

# Row 9
$ws.Range("D9").Value = "2.451.81"
$ws.Range("E9").Value = "  -5.85%  "

# Row 10
$ws.Range("E10").Value = "  -7.29%  "

# Row 11
$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  -1.86%  "

# Row 12
$ws.Range("D12").Value = "'5.37"
$ws.Range("E12").Value = "  -7.12%  "

# Row 13
$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  -6.53%  "

# Row 14
$ws.Range("D14").Value = "'25.98"
$ws.Range("E14").Value = "  -6.92%  "

# Row 15
$ws.Range("D15").Value = "2.894.99"
$ws.Range("E15").Value = "  -5.84%  "

# Row 16
$ws.Range("E16").Value = "  -8.58%  "

# Row 17
$ws.Range("D17").Value = "60.789.94"
$ws.Range("E17").Value = "  -4.15%  "

# Row 18
$ws.Range("D18").Value = "2.450.68"
$ws.Range("E18").Value = "  -5.82%  "

# Row 19
$ws.Range("E19").Value = "  -7.23%  "

# Row 20
$ws.Range("D20").Value = "'6.90"
$ws.Range("E20").Value = "  -7.49%  "

# Row 21
$ws.Range("E21").Value = "  -6.82%  "

# Row 22
$ws.Range("D22").Value = "'318.52"
$ws.Range("E22").Value = "  -6.17%  "

# Row 24
$ws.Range("D24").Value = "'63.40"
$ws.Range("E24").Value = "  -5.39%  "

# Row 25
$ws.Range("D25").Value = "'1.77"

# Row 26
$ws.Range("D26").Value = "0.0₃0967"
$ws.Range("E26").Value = "  -8.14%  "

# Row 27
$ws.Range("D27").Value = "2.573.73"
$ws.Range("E27").Value = "  -5.70%  "

# Row 29
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").Value = "'538.92"
$ws.Range("E29").Value = "  -6.53%  "

# Row 30
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.49"
$ws.Range("E30").Value = "  -3.44%  "

# Row 31
$ws.Range("D31").Value = "'8.34"
$ws.Range("E31").Value = "  -7.94%  "

# Row 32
$ws.Range("E32").Value = "  -2.29%  "

# Row 33
$ws.Range("E33").Value = "  -7.01%  "

# Row 34
$ws.Range("E34").Value = "  -7.11%  "

# Row 35
$ws.Range("E35").Value = "  -8.31%  "

# Row 36
$ws.Range("D36").Value = "'5.83"
$ws.Range("E36").Value = "  -10.34%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.84"
$ws.Range("E37").Value = "  -9.22%  "

# Row 38
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.06%  "

# Row 39
$ws.Range("E39").Value = "  -5.77%  "

# Row 40
$ws.Range("D40").Value = "'18.42"
$ws.Range("E40").Value = "  -5.88%  "

# Row 41
$ws.Range("D41").Value = "'144.91"
$ws.Range("E41").Value = "  -6.25%  "

# Row 42
$ws.Range("E42").Value = "  -0.01%  "

# Row 43
$ws.Range("D43").Value = "'1.71"
$ws.Range("E43").Value = "  -8.14%  "

# Row 44
$ws.Range("D44").Value = "'39.70"
$ws.Range("E44").Value = "  -4.47%  "

# Row 45
$ws.Range("D45").Value = "'2.29"
$ws.Range("E45").Value = "  -7.28%  "

# Row 46
$ws.Range("D46").Value = "'145.85"
$ws.Range("E46").Value = "  -7.08%  "

# Row 47
$ws.Range("D47").Value = "'3.55"
$ws.Range("E47").Value = "  -7.19%  "

# Row 48
$ws.Range("D48").Value = "'20.81"
$ws.Range("E48").Value = "  -10.82%  "

# Row 49
$ws.Range("E49").Value = "  -8.89%  "

# Row 50
$ws.Range("D50").Value = "'0.582"
$ws.Range("E50").Value = "  -6.96%  "

# Row 51
$ws.Range("D51").Value = "'0.0939"
$ws.Range("E51").Value = "  -5.40%  "
